$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 13 (shifts existing rows 13-41 down to 14-42)
$ws.Rows.Item(13).Insert()

# Fill the newly inserted row 13 with the "remove" localization entries
$ws.Range("B13").Value = "remove"
$ws.Range("C13").Value = "Remove"
$ws.Range("D13").Value = "Wissen"
$ws.Range("E13").Value = "Entfernen"

# Append the new "clickToAddNewExplanation" entry as row 43
$ws.Range("B43").Value = "clickToAddNewExplanation"
$ws.Range("C43").Value = "Click here to add explanation…"

# Update the view state to match the authored change
$ws.Range("A8").Select()
